# CSS00013 is nolonger required
# The 6th sample entry row (row 12) on the "Cora Intake Manifest" sheet
# is being removed from the intake manifest; all rows below it (the
# closing </SAMPLE ENTRIES> marker row and the "Generated by"/"Verified by"
# rows) shift up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire 12th row (the 6th sample entry: Box2 / A:6 / Container6 /
# SAMPLE_NAME6 / unique6 / subjectCode6 / test6); Excel automatically shifts
# the rows below up and adjusts the data validation ranges that reference
# the sample entry block.
$ws.Rows.Item(12).Delete()

# Update the view: make sure the sheet is active/zoomed the same as before,
# and move the selection to C26 (matching the post-edit saved state).
$ws.Activate()
$excel.ActiveWindow.Zoom = 120
$ws.Range("C26").Select()
